$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated coin list values. Price cells whose new reading is a
# plain numeric-looking string (e.g. "0.9988") are pre-formatted as Text
# (one cell at a time -- the union-range form doesn't apply to every area)
# so the literal digits (incl. trailing zeros) round-trip instead of being
# parsed into a float.
$ws.Range('D2').Value = '25.492.52'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '1.665.10'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.33'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4650'
$ws.Range('E7').Value = '  -2.52%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2577'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06140'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').Value = '1.664.16'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06944'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.64'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.370'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '74.96'
$ws.Range('E14').Value = '  +2.03%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5686'
$ws.Range('E15').Value = '  -4.67%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9999'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '25.493.89'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000006700'
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.37'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').Value = '1.877.89'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.432'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.722'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.220'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '135.81'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '14.85'
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.361'
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.701'
$ws.Range('E28').Value = '  +3.91%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '103.90'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.934'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07715'
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.607'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04298'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.618'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9450'
$ws.Range('E35').Value = '  +2.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6003'
$ws.Range('E36').Value = '  +3.38%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9296'
$ws.Range('E37').Value = '  +14.90%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.486'
$ws.Range('E38').Value = '  -2.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9993'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '103.28'
$ws.Range('E40').Value = '  +5.48%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.01461'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.819'
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3703'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.969'
$ws.Range('E44').Value = '  +6.05%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1106'
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.05253'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.126'
$ws.Range('E47').Value = '  +1.49%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '29.72'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.393'
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9983'
$ws.Range('E51').Value = '  +0.07%  '
